# "updated to flow types; conus template lib added"
#  - Rename the five "Net Growth: ..." flow types to "Net Growth Forest: ..."
#    on the "Flow Type" sheet.
#  - Update the sheet's active selection to reflect where the user ended up
#    after making the edits.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Flow Type")
$ws.Activate()

$ws.Range("A36").Value = "Net Growth Forest: Atmosphere -> Coarse Roots"
$ws.Range("A37").Value = "Net Growth Forest: Atmosphere -> Fine Roots"
$ws.Range("A38").Value = "Net Growth Forest: Atmosphere -> Foliage"
$ws.Range("A39").Value = "Net Growth Forest: Atmosphere -> Merchantable"
$ws.Range("A40").Value = "Net Growth Forest: Atmosphere -> Other Wood"

# Reflect the final cursor position / scroll state recorded in the saved
# workbook (frozen header row stays the same; selection moves to D36).
$ws.Range("D36").Select()
